$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (line_to_translate) rows 2-12: convert to Title Case
$ws.Range("B2").Value = "¿Qué Harías Si Te Besara Ahora Mismo?"
$ws.Range("B3").Value = "Für So Einen Kleinen Mann Hast Du Aber Sehr Grosse Nasenlöcher!"
$ws.Range("B4").Value = "J'Étais Assise Avec Ma Cigarette À La Main Et Je Me Suis Rendu Compte Que Je Préférerais T'Avoir Toi Dans La Main."
$ws.Range("B5").Value = "Quando Ha Visto Mia Suocera, Il Gatto Si È Arrampicato Sull'Albero E Non Ne È Ancora Sceso."
$ws.Range("B6").Value = "The Toddler’S Endless Tantrum Caused The Entire Plane Anxiety."
$ws.Range("B7").Value = "An Nomizeis Oti Qa Plhrwsw Gia Auta Tote Pou Na Deis Ti Se Perimenei Meta."
$ws.Range("B8").Value = "Als Ik Jou Mijn Telefoonnummer Geef, Zou Je Het Houden Of Weggooien?"
$ws.Range("B9").Value = "A Cama Do Meu Marido Está Cheia De Areia Preta. Porquê?"
$ws.Range("B10").Value = "It Took Him A Month To Finish The Meal."
$ws.Range("B11").Value = "¿El Coche Viajaba Rápido?"
$ws.Range("B12").Value = "Per Favore, Mi Puo' Pesare Questo Pacco?"

# Column C (translate_line) row 7 shares the same shared string as B7 (Greek sentence)
$ws.Range("C7").Value = "An Nomizeis Oti Qa Plhrwsw Gia Auta Tote Pou Na Deis Ti Se Perimenei Meta."

# Column C content rewordings (Hebrew)
$ws.Range("C6").Value = "הזעם האינסופי של הפעוט גרם לחרדה בכל המטוס."
$ws.Range("C12").Value = "בבקשה, אתה יכול לשקול את החבילה הזו בשבילי?"
